# build splits other than AndSplit #422
# Rework the "OrSplitWithAndSplit" sheet's sample data to describe a generic
# "Split" activity (TestItem_Split:0 / SplitLeft / SplitRight / EndSplit)
# instead of the old AndSplit-specific sample (TestItem_AndSplit:0 /
# AndRight1 / AndRight2), and move the active/selected tab from
# "SequenceWithOrSplit" over to "OrSplitWithAndSplit".

$wb = $excel.ActiveWorkbook

$wsSeq = $wb.Worksheets.Item("SequenceWithOrSplit")
$wsAnd = $wb.Worksheets.Item("OrSplitWithAndSplit")

# --- Rewrite the OrSplitWithAndSplit sample rows -------------------------
# Row 4: OrSplit -> AndSplit (nested split header)
$wsAnd.Cells.Item(4, 1).Value = "AndSplit"
$wsAnd.Cells.Item(4, 2).Value = ""
$wsAnd.Cells.Item(4, 3).Value = ""

# Row 5: Elementary/OrSplit/Left -> Block
$wsAnd.Cells.Item(5, 1).Value = "Block"
$wsAnd.Cells.Item(5, 2).Value = ""
$wsAnd.Cells.Item(5, 3).Value = ""

# Row 6: End -> Elementary/TestItem_Split:0/SplitLeft
$wsAnd.Cells.Item(6, 1).Value = "Elementary"
$wsAnd.Cells.Item(6, 2).Value = "TestItem_Split:0"
$wsAnd.Cells.Item(6, 3).Value = "SplitLeft"

# Row 7: AndSplit -> End
$wsAnd.Cells.Item(7, 1).Value = "End"
$wsAnd.Cells.Item(7, 2).Value = ""
$wsAnd.Cells.Item(7, 3).Value = ""

# Row 8: Block (unchanged)

# Row 9: Elementary/TestItem_AndSplit:0/AndRight1 -> Elementary/TestItem_Split:0/SplitRight
$wsAnd.Cells.Item(9, 2).Value = "TestItem_Split:0"
$wsAnd.Cells.Item(9, 3).Value = "SplitRight"

# Row 10: End (unchanged)

# Row 11: Block -> EndSplit
$wsAnd.Cells.Item(11, 1).Value = "EndSplit"
$wsAnd.Cells.Item(11, 2).Value = ""
$wsAnd.Cells.Item(11, 3).Value = ""

# Row 12: Elementary/TestItem_AndSplit:0/AndRight2 -> Block
$wsAnd.Cells.Item(12, 1).Value = "Block"
$wsAnd.Cells.Item(12, 2).Value = ""
$wsAnd.Cells.Item(12, 3).Value = ""

# Row 13: End -> Elementary/TestItem_Split:0/Right
$wsAnd.Cells.Item(13, 1).Value = "Elementary"
$wsAnd.Cells.Item(13, 2).Value = "TestItem_Split:0"
$wsAnd.Cells.Item(13, 3).Value = "Right"

# Row 14/15: End / End (unchanged)

# Column C grew a bit wider to fit "SplitRight" / "EndSplit".
$wsAnd.Columns.Item(3).ColumnWidth = 11

# --- Move the active tab / selection from SequenceWithOrSplit ------------
$wsAnd.Activate()
$wsAnd.Range("A11").Select()
